# Edit script: apply diff changes to Review_303.docx
$d = $word.ActiveDocument

# --- Paragraph 1: date + title (two runs separated by a line break) ---
$r1 = $d.Paragraphs.Item(1).Range
$ret1 = $r1.Find.Execute('🚀המאמר היומי של מייק 22.09.24: ⚡️🚀', $true, $false, $false, $false, $false, $true, 1, $false, '🚀המאמר היומי של מייק 21.09.24: ⚡️🚀', 2)
Write-Output "replace1: $ret1"

$r1b = $d.Paragraphs.Item(1).Range
$ret1b = $r1b.Find.Execute(' Quiet-STaR: Language Models Can Teach Themselves to Think Before Speaking', $true, $false, $false, $false, $false, $true, 1, $false, 'REFT: Reasoning with REinforced Fine-Tuning', 2)
Write-Output "replace1b: $ret1b"

# --- Paragraph 2: full text replacement ---
$d.Paragraphs.Item(2).Range.Text = 'ממשיכים לסקור מאמרים שסללו לכאורה נתיב ל-o1. הפעם מאמר די בסיסי יחסית שהיה שווה לסקור אותה לפני יומיים אך התעצלתי לעבור על רשימת המאמרים שבניתי כדי להבין את זה. הרווח היחיד לאלו שעוקבים אחרי סקירותיי באופן יום יומי יתבטא בכך שיהיה לכם מאוד קל להבין את הסקירה הזו אם הצלחתם להבין (בערך) את 4 הקודמות.'

# --- Paragraph 3: full text replacement ---
$d.Paragraphs.Item(3).Range.Text = 'המאמר מניח שיש בידינו דאטהסט של שאלות ושרשרת הריזונינג המובילה לתשובה (הנכונה). המאמר מציע לשפר את יכולת הריזונינג של מודל שפה בשני שלבים:'

# --- Paragraph 4: full text replacement ---
$d.Paragraphs.Item(4).Range.Text = 'אימון רגיל (Self-Supervised Fine Tuning): על כל שרשראות הריזונינג מהדאטהסט. כלומר המודל לומד לשחזר את שרשרת הריזונינג של כל שאלה ברמת הטוקן כמו ש נעשה ב-SFT הסטנדרטי.'

# --- Paragraph 5: full text replacement. The new text contains a literal
#     apostrophe ("ג'ו"), and Find/Execute's Replacement text silently
#     "smart-quotes" straight apostrophes into U+2019 - so we can't use
#     Find/Execute here. A direct Range.Text= on paragraph 5 keeps the
#     correct apostrophe, but (being a pre-existing run that had
#     xml:space="preserve" because its old text ended with a space)
#     stays marked xml:space="preserve" even though the new text needs
#     none. So: insert a brand-new paragraph after paragraph 4, set ITS
#     text (a fresh run has no stale xml:space), then delete the old
#     paragraph 5 (now shifted to index 6) entirely, mark and all. ---
$p4 = $d.Paragraphs.Item(4).Range
$p4.InsertParagraphAfter()
$newP5 = $d.Paragraphs.Item(5).Range
$newP5.Text = 'אימון של למידת פוליסי (שזה המודל עצמו) מעולם Reinforcement Learning: (מכאן נגזר שם המאמר) כאשר המודל מקבל פרס 1 אם המליח לגנרט שרשרת ריזונינג המובילה לתשובה הנכונה. תגמול צנוע הרבה יותר ניתן לתשובות מספריות לא נכונות עבור השאלות שהתשובות עליהן מספריות גם כן (כמו במאמר הקודם). תגמול 0 מתקבל בכל המקרים האחרים. אימון מתבצע עם PPO די סטנדרטי עם שערוך די סטנדרטי של פונקציית ערך V ופונקצית יתרון A (כמו במאמר המקורי של ג''ו שולמן מ-openai לשעבר)'
$oldP5 = $d.Paragraphs.Item(6).Range
$oldP5.Delete()
Write-Output "paragraph count after para5 swap: $($d.Paragraphs.Count)"

# --- Paragraph 6: entirely removed ---
$d.Paragraphs.Item(6).Range.Delete()

# --- Paragraph (now 6, was 7): URL replacement ---
$d.Paragraphs.Item(6).Range.Text = 'https://arxiv.org/pdf/2401.08967 '

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
